$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns retain text formatting so numeric-looking
# strings (e.g. "1.00") and special Unicode digits are not reinterpreted.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "72.389.69"
$ws.Range("E2").Value = "  +4.55%  "

$ws.Range("D3").Value = "4.048.00"
$ws.Range("E3").Value = "  +4.02%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "521.20"
$ws.Range("E5").Value = "  -0.99%  "

$ws.Range("D6").Value = "147.24"
$ws.Range("E6").Value = "  +3.02%  "

$ws.Range("D7").Value = "0.623"
$ws.Range("E7").Value = "  +2.28%  "

$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "0.737"
$ws.Range("E9").Value = "  +2.56%  "

$ws.Range("D10").Value = "0.175"
$ws.Range("E10").Value = "  +2.52%  "

$ws.Range("D11").Value = "0.0000334"
$ws.Range("E11").Value = "  +0.90%  "

$ws.Range("D12").Value = "47.62"
$ws.Range("E12").Value = "  +13.56%  "

$ws.Range("D13").Value = "10.89"
$ws.Range("E13").Value = "  +6.90%  "

$ws.Range("D14").Value = "4.678.90"
$ws.Range("E14").Value = "  +3.45%  "

$ws.Range("D15").Value = "4.041.77"
$ws.Range("E15").Value = "  -1.29%  "

$ws.Range("D16").Value = "21.18"
$ws.Range("E16").Value = "  +7.76%  "

$ws.Range("D17").Value = "14.16"
$ws.Range("E17").Value = "  +2.91%  "

$ws.Range("E18").Value = "  -1.71%  "

$ws.Range("E19").Value = "  -2.16%  "

$ws.Range("D20").Value = "72.310.92"
$ws.Range("E20").Value = "  +4.45%  "

$ws.Range("D21").Value = "437.91"
$ws.Range("E21").Value = "  +3.17%  "

$ws.Range("D22").Value = "97.77"
$ws.Range("E22").Value = "  +11.56%  "

$ws.Range("D23").Value = "3.54"
$ws.Range("E23").Value = "  +5.98%  "

$ws.Range("D24").Value = "14.71"
$ws.Range("E24").Value = "  +4.06%  "

$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "4.02"
$ws.Range("E25").Value = "  -1.17%  "

$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").Value = "11.93"
$ws.Range("E26").Value = "  +2.53%  "

$ws.Range("D27").Value = "11.29"
$ws.Range("E27").Value = "  +7.28%  "

$ws.Range("D28").Value = "37.08"
$ws.Range("E28").Value = "  +3.17%  "

$ws.Range("D29").Value = "3.09"
$ws.Range("E29").Value = "  +10.02%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "13.53"
$ws.Range("E30").Value = "  +3.21%  "

$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "692.61"
$ws.Range("E31").Value = "  -0.26%  "

$ws.Range("E32").Value = "  +2.81%  "

$ws.Range("D33").Value = "7.01"

$ws.Range("D34").Value = "68.09"
$ws.Range("E34").Value = "  +0.61%  "

$ws.Range("E35").Value = "  +6.89%  "

$ws.Range("D36").Value = "0.439"
$ws.Range("E36").Value = "  -0.72%  "

$ws.Range("D37").Value = "3.65"
$ws.Range("E37").Value = "  +23.84%  "

$ws.Range("D38").Value = "40.68"
$ws.Range("E38").Value = "  +1.14%  "

$ws.Range("E39").Value = "  +3.06%  "

$ws.Range("E40").Value = "  +0.30%  "

$ws.Range("D41").Value = "0.997"
$ws.Range("E41").Value = "  -0.41%  "

$ws.Range("D42").Value = "0.0488"
$ws.Range("E42").Value = "  +1.76%  "

$ws.Range("D43").Value = "3.16"
$ws.Range("E43").Value = "  +5.13%  "

$ws.Range("D44").Value = "2.77"
$ws.Range("E44").Value = "  -1.21%  "

$ws.Range("D45").Value = "3.51"
$ws.Range("E45").Value = "  +6.00%  "

$ws.Range("D46").Value = "0.146"
$ws.Range("E46").Value = "  +4.40%  "

$ws.Range("E47").Value = "  +2.10%  "

$ws.Range("D48").Value = "9.06"
$ws.Range("E48").Value = "  +8.45%  "

$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").Value = "0.000273"
$ws.Range("E49").Value = "  +21.04%  "

$ws.Range("B50").Value = "LidoDAOToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D50").Value = "3.32"
$ws.Range("E50").Value = "  +1.48%  "

$ws.Range("D51").Value = "0.0₆0339"
$ws.Range("E51").Value = "  -0.66%  "

